$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.853.85'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.527.22'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.98%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.26%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.524.74'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.147'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.85%  '

$ws.Range("E11").Value = '  -1.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.98'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.03%  '

$ws.Range("E13").Value = '  +0.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.986.28'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.16'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.755.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.24%  '

$ws.Range("E17").Value = '  +0.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.518.72'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.10%  '

$ws.Range("E19").Value = '  +1.53%  '

$ws.Range("E20").Value = '  +0.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '353.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.18'
$ws.Range("D22").Style = "Normal"

$ws.Range("E23").Value = '  +0.03%  '

$ws.Range("E24").Value = '  +1.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.23'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.42%  '

$ws.Range("E26").Value = '  -5.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.74%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.655.77'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.996'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '521.52'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0892'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.50%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.81'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.30%  '

$ws.Range("E33").Value = '  +0.53%  '

$ws.Range("E34").Value = '  +0.90%  '

$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("E36").Value = '  +0.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.01'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("E38").Value = '  +1.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.70'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.79'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.67%  '

$ws.Range("E41").Value = '  -1.19%  '

$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.85'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.325'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.94%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.48%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '153.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.14%  '

$ws.Range("E48").Value = '  +2.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0256'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.20%  '

$ws.Range("E50").Value = '  +1.80%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0741'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.28%  '
